$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-6 in place with the new fixture data (column B/Date is
# untouched on every row, so it is skipped to avoid Excel's automatic
# date literal -> serial-number conversion on write).

# Row 3
$ws.Cells.Item(3,1).Value = 'Brazilian Serie A'
$ws.Cells.Item(3,3).Value = '19:30:00'
$ws.Cells.Item(3,4).Value = 'Vasco da Gama'
$ws.Cells.Item(3,5).Value = 'Internacional'
$ws.Cells.Item(3,8).Value = 1000
$ws.Cells.Item(3,10).Value = 1000
$ws.Cells.Item(3,18).Value = 0
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0
$ws.Cells.Item(3,23).Value = 500
$ws.Cells.Item(3,36).Value = 1000
$ws.Cells.Item(3,37).Value = 1000
$ws.Cells.Item(3,38).Value = 1000
$ws.Cells.Item(3,40).Value = 1.01

# Row 4
$ws.Cells.Item(4,3).Value = '21:00:00'
$ws.Cells.Item(4,4).Value = 'Academia de Balompie Boliviano'
$ws.Cells.Item(4,5).Value = 'Always Ready'
$ws.Cells.Item(4,6).Value = 100
$ws.Cells.Item(4,7).Value = 1000
$ws.Cells.Item(4,8).Value = 1.01
$ws.Cells.Item(4,9).Value = 1.01
$ws.Cells.Item(4,10).Value = 100
$ws.Cells.Item(4,11).Value = 1000
$ws.Cells.Item(4,18).Value = 0
$ws.Cells.Item(4,19).Value = 0
$ws.Cells.Item(4,20).Value = 1.37
$ws.Cells.Item(4,21).Value = 1.01
$ws.Cells.Item(4,22).Value = 1.01
$ws.Cells.Item(4,23).Value = 1.01
$ws.Cells.Item(4,40).Value = 1000
$ws.Cells.Item(4,41).Value = 1.9

# Row 5
$ws.Cells.Item(5,1).Value = 'Ecuadorian Serie A'
$ws.Cells.Item(5,3).Value = '21:00:00'
$ws.Cells.Item(5,4).Value = 'LDU'
$ws.Cells.Item(5,5).Value = 'Independiente (Ecu)'
$ws.Cells.Item(5,6).Value = 5.8
$ws.Cells.Item(5,7).Value = 11.5
$ws.Cells.Item(5,8).Value = 12
$ws.Cells.Item(5,9).Value = 15.5
$ws.Cells.Item(5,10).Value = 1.27
$ws.Cells.Item(5,11).Value = 1.3
$ws.Cells.Item(5,16).Value = 1.3
$ws.Cells.Item(5,17).Value = 3.3
$ws.Cells.Item(5,18).Value = 1.01
$ws.Cells.Item(5,19).Value = 36
$ws.Cells.Item(5,22).Value = 1.08
$ws.Cells.Item(5,23).Value = 1.16
$ws.Cells.Item(5,29).Value = 1.3
$ws.Cells.Item(5,30).Value = 28
$ws.Cells.Item(5,31).Value = 1000
$ws.Cells.Item(5,33).Value = 1000
$ws.Cells.Item(5,34).Value = 520
$ws.Cells.Item(5,35).Value = 1000
$ws.Cells.Item(5,37).Value = 1000
$ws.Cells.Item(5,38).Value = 1000

# Row 6
$ws.Cells.Item(6,3).Value = '21:30:00'
$ws.Cells.Item(6,4).Value = 'Santos'
$ws.Cells.Item(6,5).Value = 'Sport Recife'
$ws.Cells.Item(6,6).Value = 1.02
$ws.Cells.Item(6,7).Value = 1.03
$ws.Cells.Item(6,8).Value = 240
$ws.Cells.Item(6,9).Value = 260
$ws.Cells.Item(6,10).Value = 38
$ws.Cells.Item(6,11).Value = 40
$ws.Cells.Item(6,16).Value = 5
$ws.Cells.Item(6,17).Value = 1.24
$ws.Cells.Item(6,18).Value = 1.79
$ws.Cells.Item(6,19).Value = 2.2
$ws.Cells.Item(6,20).Value = 2.66
$ws.Cells.Item(6,21).Value = 1.56
$ws.Cells.Item(6,23).Value = 34
$ws.Cells.Item(6,32).Value = 5.3
$ws.Cells.Item(6,33).Value = 10
$ws.Cells.Item(6,34).Value = 44
$ws.Cells.Item(6,35).Value = 360
$ws.Cells.Item(6,36).Value = 4.3
$ws.Cells.Item(6,37).Value = 9.199999999999999
$ws.Cells.Item(6,38).Value = 42
$ws.Cells.Item(6,39).Value = 320
$ws.Cells.Item(6,40).Value = 3.55

# Remove the trailing three fixtures (old rows 7-9); the sheet now ends at row 6.
$ws.Range("A7:A9").EntireRow.Delete()
